$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.987.93"
$ws.Range("E2").Value = '  +2.69%  '
$ws.Range("D3").Value = "'2.054.24"
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'229.89"
$ws.Range("E5").Value = '  +1.98%  '
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("D7").Value = "'58.73"
$ws.Range("E7").Value = '  +7.69%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +3.35%  '
$ws.Range("D10").Value = "'0.0813"
$ws.Range("E10").Value = '  +4.71%  '
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = '  +2.14%  '
$ws.Range("D12").Value = "'2.358.05"
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("E13").Value = '  +4.70%  '
$ws.Range("E14").Value = '  +5.57%  '
$ws.Range("E15").Value = '  +2.63%  '
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = "'2.045.22"
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = "'37.909.16"
$ws.Range("E18").Value = '  +2.84%  '
$ws.Range("D19").Value = "'6.36"
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = "'69.76"
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("E21").Value = '  +3.42%  '
$ws.Range("D22").Value = "'224.57"
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = "'2.42"
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("E25").Value = '  +4.16%  '
$ws.Range("D26").Value = "'166.60"
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("E28").Value = '  +5.52%  '
$ws.Range("E29").Value = '  +2.89%  '
$ws.Range("E30").Value = '  +2.96%  '
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("D32").Value = "'4.53"
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("E33").Value = '  +3.87%  '
$ws.Range("E34").Value = '  +10.81%  '
$ws.Range("D35").Value = "'0.0608"
$ws.Range("E35").Value = '  +1.67%  '
$ws.Range("D36").Value = "'2.33"
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("E37").Value = '  +13.26%  '
$ws.Range("E38").Value = '  +6.29%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = "'1.520.66"
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = '  +4.47%  '
$ws.Range("D43").Value = "'97.07"
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("D44").Value = "'16.53"
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").Value = "'0.0921"
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = "'4.09"
$ws.Range("E47").Value = '  +16.33%  '
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("D50").Value = "'7.09"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = "'2.246.38"
$ws.Range("E51").Value = '  +2.35%  '
